$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, borders, centered) from an existing header
# cell (H1) onto the two new header cells so I1/J1 match the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Header labels for the new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New column data
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 2

$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 8

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 5

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 5

$ws.Range("I6").Value = 4
$ws.Range("J6").Value = 5
